$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Language most often" wording fixes ---------------------------------
$ws.Range("A18").Value = "English spoken most often at home"
$ws.Range("A19").Value = "French spoken most often at home"
$ws.Range("B20").Value = "Non-official language spoken at home includes people who provided a single-response of the language spoken most often at home, which is not English or French. Shown are the top ten (if relevant) languages most often spoken at home."
$ws.Range("B22").Value = "Median household income is the median total income for private households in 2015. If relevant, the median is estimated using bucket counts and assuming a uniform distribution."

# --- View/selection state (scrolled down, new selection) -----------------
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C27").Select()
